$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83").Value = "Zwane deserved a fair chance as Chiefs coach – Mabedi"
$ws.Range("B83").Value = "Zwane"
$ws.Range("C83").Value = "Q4800783"
$ws.Range("D83").Value = "South African footballer"

$ws.Range("A84").Value = "President Mnangagwa officiates 4th ZNDU graduation ceremony"
$ws.Range("B84").Value = "President Mnangagwa"
$ws.Range("C84").Value = "Q510523"
$ws.Range("D84").Value = "President of the Republic of Zimbabwe"
